$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (shared string change)
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 21:22"

# Country data table updates: reordered countries (new data pushed up, displaced
# countries keep their old figures but shift down a row) plus independent refreshes
# for a handful of countries whose totals were updated in this snapshot.

$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 388757
$ws.Cells.Item(4, 3).Value = 21753
$ws.Cells.Item(4, 4).Value = 21431
$ws.Cells.Item(4, 5).Value = 354926
$ws.Cells.Item(4, 6).Value = 9165
$ws.Cells.Item(4, 7).Value = 1529
$ws.Cells.Item(4, 8).Value = 12400

$ws.Cells.Item(17, 1).Value = "Austria"
$ws.Cells.Item(17, 2).Value = 12633
$ws.Cells.Item(17, 3).Value = 336
$ws.Cells.Item(17, 4).Value = 4046
$ws.Cells.Item(17, 5).Value = 8344
$ws.Cells.Item(17, 6).Value = 243
$ws.Cells.Item(17, 7).Value = 23
$ws.Cells.Item(17, 8).Value = 243

$ws.Cells.Item(18, 1).Value = "Brasil"
$ws.Cells.Item(18, 2).Value = 12632
$ws.Cells.Item(18, 3).Value = 449
$ws.Cells.Item(18, 4).Value = 127
$ws.Cells.Item(18, 5).Value = 11917
$ws.Cells.Item(18, 6).Value = 296
$ws.Cells.Item(18, 7).Value = 24
$ws.Cells.Item(18, 8).Value = 588

$ws.Cells.Item(39, 1).Value = "Peru"
$ws.Cells.Item(39, 2).Value = 2954
$ws.Cells.Item(39, 3).Value = 393
$ws.Cells.Item(39, 4).Value = 1301
$ws.Cells.Item(39, 5).Value = 1546
$ws.Cells.Item(39, 6).Value = 109
$ws.Cells.Item(39, 7).Value = 15
$ws.Cells.Item(39, 8).Value = 107

$ws.Cells.Item(58, 1).Value = "Egipto"
$ws.Cells.Item(58, 2).Value = 1450
$ws.Cells.Item(58, 3).Value = 128
$ws.Cells.Item(58, 4).Value = 276
$ws.Cells.Item(58, 5).Value = 1080
$ws.Cells.Item(58, 6).Value = 0
$ws.Cells.Item(58, 7).Value = 9
$ws.Cells.Item(58, 8).Value = 94

$ws.Cells.Item(87, 1).Value = "Costa Rica"
$ws.Cells.Item(87, 2).Value = 483
$ws.Cells.Item(87, 3).Value = 16
$ws.Cells.Item(87, 4).Value = 24
$ws.Cells.Item(87, 5).Value = 457
$ws.Cells.Item(87, 6).Value = 14
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 2

$ws.Cells.Item(143, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(143, 2).Value = 47
$ws.Cells.Item(143, 3).Value = 5
$ws.Cells.Item(143, 4).Value = 0
$ws.Cells.Item(143, 5).Value = 47
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 0

$ws.Cells.Item(144, 1).Value = "Congo"
$ws.Cells.Item(144, 2).Value = 45
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 2
$ws.Cells.Item(144, 5).Value = 38
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 5

$ws.Cells.Item(145, 1).Value = "Macao"
$ws.Cells.Item(145, 2).Value = 44
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 10
$ws.Cells.Item(145, 5).Value = 34
$ws.Cells.Item(145, 6).Value = 1
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 0

$ws.Cells.Item(167, 1).Value = "Angola"
$ws.Cells.Item(167, 2).Value = 17
$ws.Cells.Item(167, 3).Value = 1
$ws.Cells.Item(167, 4).Value = 2
$ws.Cells.Item(167, 5).Value = 13
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 2

$ws.Cells.Item(168, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(168, 2).Value = 16
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 3
$ws.Cells.Item(168, 5).Value = 13
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 0

$ws.Cells.Item(169, 1).Value = "Namibia"
$ws.Cells.Item(169, 2).Value = 16
$ws.Cells.Item(169, 3).Value = 0
$ws.Cells.Item(169, 4).Value = 3
$ws.Cells.Item(169, 5).Value = 13
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = 0

$ws.Cells.Item(180, 1).Value = "Seychelles"
$ws.Cells.Item(180, 2).Value = 11
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 5).Value = 11
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0

$ws.Cells.Item(181, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(181, 2).Value = 11
$ws.Cells.Item(181, 3).Value = 1
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 11
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 0

$ws.Cells.Item(184, 1).Value = "Mozambique"
$ws.Cells.Item(184, 2).Value = 10
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 1
$ws.Cells.Item(184, 5).Value = 9
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

$ws.Cells.Item(185, 1).Value = "Surinam"
$ws.Cells.Item(185, 2).Value = 10
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 5).Value = 9
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 1
